$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diary is now kept by Ágúst Heiðar Hannesson instead of Hildur Sif
# Thorarensen, so swap her name/email for his in the "Vika 8" sheet header.
$ws.Range("C2").Value = "Ágúst Heiðar Hannesson"
$ws.Range("G2").Value = "ahh19@hi.is"

# The "Vika 8" summary-table label lower on the sheet is unaffected by the
# name change but moves in the shared-string table as a side effect.
$ws.Range("B44").Value = "Vika 8"

# Vika 8 now has some logged hours under "Rannsóknir" (E27) and "Forritun"
# (E29); the SUM formulas in column J (and the D57/D59/D62 weekly totals)
# recalculate automatically.
$ws.Range("E27").Value = 120
$ws.Range("E29").Value = 180

# Restore the view to the top of the sheet with C3 selected instead of the
# previous scrolled-down E29 selection.
$ws.Range("C3").Select()
